$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New test case row (TCID ENW043) appended after the existing last row (43).
# Shared-string insertion order matters (matches authoring order in the diff):
#   1) TCID            -> column A
#   2) Description      -> column C
#   3) Jira id           -> column B
$ws.Cells.Item(44, 1).Value = "ENW043"
$ws.Cells.Item(44, 3).Value = "Verify that the user signed in to community enabled version of Endnote and having invalid Neon session will be taken to Privacy page seamlessly by clicking on the Privacy in profile flyout."
$ws.Cells.Item(44, 2).Value = "`nOPQA-3617"
$ws.Cells.Item(44, 4).Value = "Y"

# Match the formatting used by the rest of the table: thin border all around,
# wrapped text for the Jira id / Description columns, row height 45.
$ws.Range("A44:E44").Borders.LineStyle = 1
$ws.Range("B44:C44").WrapText = $true
$ws.Rows.Item(44).RowHeight = 45

# Mirror the cursor move captured in the workbook (selection ends up on D49).
[void]$ws.Range("D49").Select()
